$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Raef%Fadel%NULL%0,         Austin R%Morrison%NULL%1,         Amit%Vahia%NULL%0,         Zachary R%Smith%NULL%1,         Zohra%Chaudhry%NULL%0,         Pallavi%Bhargava%NULL%1,         Joseph%Miller%NULL%0,         Rachel M%Kenney%NULL%1,         George%Alangaden%NULL%0,         Mayur S%Ramesh%mramesh1@hfhs.org%1,         Varidhi%Nauriyal%NULL%1,         Jayanth%Lakshmikanth%NULL%1,         Asif%Abdul Hamed%NULL%1,         Owais%Nadeem%NULL%1,         Kristin%Griebe%NULL%1,         Joseph M%Johnson%NULL%1,         Patrick%Bradley%NULL%1,         Junior%Uduman%NULL%1,         Sara%Hegab%NULL%1,         Jennifer%Swiderek%NULL%1,         Amanda%Godfrey%NULL%1,         Jeffrey%Jennings%NULL%1,         Jayna%Gardner-Gray%NULL%1,         Adam%Ackerman%NULL%1,         Jonathan%Lezotte%NULL%1,         Joseph%Ruhala%NULL%1,         Linoj%Samuel%NULL%1,         Robert J%Tibbetts%NULL%1,         Indira%Brar%NULL%0,         John%McKinnon%NULL%1,         Geehan%Suleyman%NULL%1,         Nicholas%Yared%NULL%1,         Erica%Herc%NULL%1,         Jonathan%Williams%NULL%1,         Odaliz Abreu%Lanfranco%NULL%1,         Anne%Chen%NULL%0,         Marcus%Zervos%NULL%1,         Eric%Scher%NULL%1]'
$ws.Range("E3").Value = '[Xiaowei%Fang%NULL%2,         Qing%Mei%NULL%2,         Tianjun%Yang%NULL%2,         Lei%Li%NULL%0,         Yinzhong%Wang%NULL%2,         Fei%Tong%NULL%2,         Shike%Geng%NULL%2,         Aijun%Pan%NULL%3]'
$ws.Range("E4").Value = '[Xiaowei%Fang%NULL%0,         Qing%Mei%NULL%0,         Tianjun%Yang%NULL%0,         Lei%Li%NULL%0,         Yinzhong%Wang%NULL%0,         Fei%Tong%NULL%0,         Shike%Geng%NULL%0,         Aijun%Pan%NULL%0]'
$ws.Range("E5").Value = '[Yun%Ling%NULL%2,         Shui-Bao%Xu%NULL%2,         Yi-Xiao%Lin%NULL%2,         Di%Tian%NULL%2,         Zhao-Qin%Zhu%NULL%2,         Fa-Hui%Dai%NULL%2,         Fan%Wu%NULL%2,         Zhi-Gang%Song%NULL%2,         Wei%Huang%NULL%2,         Jun%Chen%NULL%0,         Bi-Jie%Hu%NULL%2,         Sheng%Wang%NULL%2,         En-Qiang%Mao%NULL%2,         Lei%Zhu%NULL%2,         Wen-Hong%Zhang%NULL%2,         Hong-Zhou%Lu%NULL%2,         Peng%Lyu%NULL%4,         Peng%Lyu%NULL%0]'
$ws.Range("E6").Value = '[Xiaofan%Lu%NULL%0,         Taige%Chen%NULL%1,         Yang%Wang%NULL%1,         Jun%Wang%NULL%2,         Fangrong%Yan%f.r.yan@163.com%1]'
$ws.Range("E7").Value = '[Yin%Wang%NULL%0,         Weiwei%Jiang%NULL%1,         Qi%He%NULL%1,         Cheng%Wang%NULL%3,         Baoju%Wang%NULL%1,         Pan%Zhou%NULL%1,         Nianguo%Dong%dongnianguo63@gmail.com%1,         Qiaoxia%Tong%2013xh0859@hust.edu.cn%1]'
$ws.Range("E8").Value = '[ J.%Wu%null%1,          J.% Huang%null%1,          G.% Zhu%null%1,          Y.% Liu%null%2,          H.% Xiao%null%1,          Q.% Zhou%null%1,          X.% Si%null%1,          H.% Yi%null%1,          C.% Wang%null%0,          D.% Yang%null%1,          S.% Chen%null%1,          X.% Liu%null%1,          Z.% Liu%null%1,          Q.% Wang%null%1,          Q.% Lv%null%1,          Y.% Huang%null%1,          Y.% Yu%null%1,          X.% Guan%null%1,          Y.% Li%null%3,          K.% Nirantharakumar%null%1,          K.% Cheng%null%1,          S.% Peng%null%1,          H. % Xiao%null%1]'
$ws.Range("E10").Value = '[Xiaobo%Yang%NULL%0,         Yuan%Yu%NULL%0,         Jiqian%Xu%NULL%0,         Huaqing%Shu%NULL%0,         Jia''an%Xia%NULL%0,         Hong%Liu%NULL%0,         Yongran%Wu%NULL%0,         Lu%Zhang%NULL%0,         Zhui%Yu%NULL%0,         Minghao%Fang%NULL%0,         Ting%Yu%NULL%0,         Yaxin%Wang%NULL%0,         Shangwen%Pan%NULL%0,         Xiaojing%Zou%NULL%0,         Shiying%Yuan%NULL%0,         You%Shang%NULL%0]'
$ws.Range("E11").Value = '[Lei%Zha%NULL%0,         Shirong%Li%NULL%1,         Lingling%Pan%NULL%1,         Boris%Tefsen%NULL%1,         Yeshan%Li%NULL%1,         Neil%French%NULL%1,         Liyun%Chen%NULL%1,         Gang%Yang%NULL%1,         Elmer V%Villanueva%Villanueva@xjtlu.edu.cn%1]'
$ws.Range("E12").Value = '[Wei%Zhou%NULL%2,         Yisi%Liu%NULL%2,         Dongdong%Tian%NULL%2,         Cheng%Wang%NULL%0,         Sa%Wang%NULL%2,         Jing%Cheng%NULL%2,         Ming%Hu%huming74@163.com%0,         Minghao%Fang%fangmh@tjh.tjmu.edu.cn%0,         Yue%Gao%gaoyue@bmi.ac.cn%3]'
$ws.Range("E13").Value = '[Wei%Zhou%NULL%0,         Yisi%Liu%NULL%0,         Dongdong%Tian%NULL%0,         Cheng%Wang%NULL%0,         Sa%Wang%NULL%0,         Jing%Cheng%NULL%0,         Ming%Hu%huming74@163.com%0,         Minghao%Fang%fangmh@tjh.tjmu.edu.cn%0,         Yue%Gao%gaoyue@bmi.ac.cn%0]'
$ws.Range("E14").Value = '[Rong-chang%Chen%NULL%1,         Xiao-ping%Tang%NULL%1,         Shou-yong%Tan%NULL%1,         Bi-ling%Liang%NULL%1,         Zhuo-yue%Wan%NULL%1,         Ji-qian%Fang%NULL%1,         Nanshan%Zhong%NULL%0]'
$ws.Range("E15").Value = '[Justin%Coyle%NULL%1,         Efehi%Igbinomwanhia%NULL%1,         Alejandro%Sanchez-Nadales%alejandro.sanchez2@advocatehealth.com%1,         Sorin%Danciu%NULL%1,         Chae%Chu%NULL%1,         Nishit%Shah%NULL%1]'
$ws.Range("E16").Value = '[Yuanyuan%Han%NULL%1,         Mao%Jiang%NULL%1,         Da%Xia%NULL%1,         Lichao%He%NULL%1,         Xin%Lv%NULL%1,         Xiaohua%Liao%NULL%1,         Jie%Meng%mengjie@csu.edu.cn%1]'
$ws.Range("E17").Value = '[Kui%Liu%NULL%0,         Yuan-Yuan%Fang%NULL%2,         Yan%Deng%NULL%2,         Wei%Liu%NULL%2,         Mei-Fang%Wang%NULL%1,         Jing-Ping%Ma%NULL%1,         Wei%Xiao%NULL%2,         Ying-Nan%Wang%NULL%1,         Min-Hua%Zhong%NULL%1,         Cheng-Hong%Li%NULL%1,         Guang-Cai%Li%NULL%1,         Hui-Guo%Liu%NULL%3,         Xiu-Yuan%Hao%NULL%4,         Xiu-Yuan%Hao%NULL%0,         Pei-Fang%Wei%NULL%0]'
